$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 174 (01-05-2021)
$ws.Range("B174").Value = 62522
$ws.Range("H174").Value = 6651
$ws.Range("I174").Value = 28433

# Row 175 (01-06-2021)
$ws.Range("B175").Value = 59078
$ws.Range("H175").Value = 6351
$ws.Range("I175").Value = 27930

# Row 176 (01-07-2021)
$ws.Range("B176").Value = 56397
$ws.Range("H176").Value = 6905
$ws.Range("I176").Value = 27456

# Row 177 (01-08-2021) - update existing value and add missing columns
$ws.Range("B177").Value = 56850
$ws.Range("C177").Value = 264
$ws.Range("D177").Value = 3561
$ws.Range("E177").Value = 8979
$ws.Range("F177").Value = 3595
$ws.Range("G177").Value = 6496
$ws.Range("H177").Value = 6584
$ws.Range("I177").Value = 27371
